$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table (B2:E51) with refreshed data.
# NumberFormat="@" forces the written Value2 to be stored as text (matching the
# source workbook convention where prices/percentages are plain strings, not
# numeric values), and ClearFormats() afterwards drops the temporary number
# format again so the cell keeps its original (default) style.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "26.742.52"
Set-TextValue $ws.Range("E2") "  +0.62%  "
Set-TextValue $ws.Range("D3") "1.648.14"
Set-TextValue $ws.Range("E3") "  +1.13%  "
Set-TextValue $ws.Range("E4") "  +0.12%  "
Set-TextValue $ws.Range("D5") "215.42"
Set-TextValue $ws.Range("E5") "  +1.22%  "
Set-TextValue $ws.Range("D6") "0.503"
Set-TextValue $ws.Range("E6") "  +1.63%  "
Set-TextValue $ws.Range("E7") "  +0.13%  "
Set-TextValue $ws.Range("E8") "  -0.37%  "
Set-TextValue $ws.Range("D9") "0.0627"
Set-TextValue $ws.Range("E9") "  +0.79%  "
Set-TextValue $ws.Range("D10") "19.23"
Set-TextValue $ws.Range("E10") "  +1.55%  "
Set-TextValue $ws.Range("D11") "0.0845"
Set-TextValue $ws.Range("E11") "  +1.02%  "
Set-TextValue $ws.Range("D12") "1.877.46"
Set-TextValue $ws.Range("E12") "  +1.09%  "
Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "4.18"
Set-TextValue $ws.Range("E13") "  +2.72%  "
Set-TextValue $ws.Range("B14") "WrappedEther"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D14") "1.613.49"
Set-TextValue $ws.Range("E14") "  -1.11%  "
Set-TextValue $ws.Range("D15") "0.531"
Set-TextValue $ws.Range("E15") "  +1.38%  "
Set-TextValue $ws.Range("D16") "65.96"
Set-TextValue $ws.Range("E16") "  +4.77%  "
Set-TextValue $ws.Range("D17") "26.777.91"
Set-TextValue $ws.Range("E17") "  +0.80%  "
Set-TextValue $ws.Range("D18") "0.0₃0750"
Set-TextValue $ws.Range("E18") "  +1.45%  "
Set-TextValue $ws.Range("D19") "219.07"
Set-TextValue $ws.Range("E19") "  +4.73%  "
Set-TextValue $ws.Range("E20") "  +0.20%  "
Set-TextValue $ws.Range("E21") "  +2.13%  "
Set-TextValue $ws.Range("E22") "  +3.16%  "
Set-TextValue $ws.Range("E23") "  +0.19%  "
Set-TextValue $ws.Range("D24") "2.10"
Set-TextValue $ws.Range("E24") "  +9.93%  "
Set-TextValue $ws.Range("D25") "148.07"
Set-TextValue $ws.Range("E25") "  +0.19%  "
Set-TextValue $ws.Range("E26") "  +0.22%  "
Set-TextValue $ws.Range("E27") "  -0.91%  "
Set-TextValue $ws.Range("D28") "7.00"
Set-TextValue $ws.Range("E28") "  +2.15%  "
Set-TextValue $ws.Range("D29") "15.78"
Set-TextValue $ws.Range("E29") "  +2.73%  "
Set-TextValue $ws.Range("D30") "0.0520"
Set-TextValue $ws.Range("E30") "  +0.87%  "
Set-TextValue $ws.Range("D31") "1.18"
Set-TextValue $ws.Range("E31") "  +0.67%  "
Set-TextValue $ws.Range("D32") "3.40"
Set-TextValue $ws.Range("E32") "  +4.75%  "
Set-TextValue $ws.Range("D33") "3.02"
Set-TextValue $ws.Range("E33") "  +3.85%  "
Set-TextValue $ws.Range("D34") "1.277.29"
Set-TextValue $ws.Range("E34") "  +9.78%  "
Set-TextValue $ws.Range("E35") "  +3.14%  "
Set-TextValue $ws.Range("D36") "0.0181"
Set-TextValue $ws.Range("E36") "  +4.64%  "
Set-TextValue $ws.Range("D37") "2.40"
Set-TextValue $ws.Range("E37") "  +1.35%  "
Set-TextValue $ws.Range("D38") "0.814"
Set-TextValue $ws.Range("E38") "  +1.30%  "
Set-TextValue $ws.Range("D39") "0.519"
Set-TextValue $ws.Range("E39") "  +3.03%  "
Set-TextValue $ws.Range("E40") "  +0.19%  "
Set-TextValue $ws.Range("D41") "2.29"
Set-TextValue $ws.Range("E41") "  -1.44%  "
Set-TextValue $ws.Range("E42") "  +2.27%  "
Set-TextValue $ws.Range("D43") "5.40"
Set-TextValue $ws.Range("E43") "  +0.51%  "
Set-TextValue $ws.Range("D44") "1.786.63"
Set-TextValue $ws.Range("E44") "  +1.27%  "
Set-TextValue $ws.Range("D45") "93.37"
Set-TextValue $ws.Range("E45") "  +0.89%  "
Set-TextValue $ws.Range("D46") "1.62"
Set-TextValue $ws.Range("E46") "  +5.74%  "
Set-TextValue $ws.Range("D47") "56.21"
Set-TextValue $ws.Range("E47") "  +3.31%  "
Set-TextValue $ws.Range("B48") "Cronos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.0516"
Set-TextValue $ws.Range("E48") "  +0.70%  "
Set-TextValue $ws.Range("B49") "EnergySwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "7.73"
Set-TextValue $ws.Range("E49") "  +1.77%  "
Set-TextValue $ws.Range("B50") "Algorand"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D50") "0.0977"
Set-TextValue $ws.Range("E50") "  +3.96%  "
Set-TextValue $ws.Range("B51") "Mantle"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D51") "0.408"
Set-TextValue $ws.Range("E51") "  -0.15%  "
